{"js": "const replacements = [\n  [\"2023-11-29 Wednesday\", \"2023-11-30 Thursday\"],\n  [\"91\u00d753=4823\", \"58\u00d760=3480\"],\n  [\"60\u00d758=3480\", \"13\u00d738=494\"],\n  [\"36\u00d749=1764\", \"18\u00d766=1188\"],\n  [\"29\u00d736=1044\", \"57\u00d737=2109\"],\n  [\"67\u00d756=3752\", \"92\u00d793=8556\"],\n  [\"85\u00d740=3400\", \"39\u00d789=3471\"],\n  [\"83\u00d771=5893\", \"85\u00d772=6120\"],\n  [\"49\u00d728=1372\", \"21\u00d746=966\"],\n  [\"86\u00d749=4214\", \"48\u00d760=2880\"],\n  [\"78\u00d777=6006\", \"18\u00d755=990\"],\n  [\"67\u00d778=5226\", \"60\u00d777=4620\"],\n  [\"93\u00d756=5208\", \"44\u00d777=3388\"],\n  [\"62\u00d741=2542\", \"63\u00d723=1449\"],\n  [\"61\u00d750=3050\", \"68\u00d795=6460\"],\n  [\"56\u00d799=5544\", \"75\u00d736=2700\"],\n  [\"14\u00d785=1190\", \"33\u00d759=1947\"],\n  [\"97\u00d712=1164\", \"22\u00d796=2112\"],\n  [\"25\u00d728=700\", \"83\u00d717=1411\"],\n  [\"48\u00d731=1488\", \"84\u00d734=2856\"],\n  [\"14\u00d753=742\", \"14\u00d761=854\"],\n  [\"62\u00d746=2852\", \"82\u00d752=4264\"],\n  [\"91\u00d758=5278\", \"97\u00d716=1552\"],\n  [\"17\u00d785=1445\", \"65\u00d762=4030\"],\n  [\"57\u00d794=5358\", \"33\u00d742=1386\"],\n  [\"19\u00d776=1444\", \"17\u00d763=1071\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n\nfunction Replace-Text($oldText, $newText) {\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nReplace-Text \"2023-11-29 Wednesday\" \"2023-11-30 Thursday\"\nReplace-Text \"91\u00d753=4823\" \"58\u00d760=3480\"\nReplace-Text \"60\u00d758=3480\" \"13\u00d738=494\"\nReplace-Text \"36\u00d749=1764\" \"18\u00d766=1188\"\nReplace-Text \"29\u00d736=1044\" \"57\u00d737=2109\"\nReplace-Text \"67\u00d756=3752\" \"92\u00d793=8556\"\nReplace-Text \"85\u00d740=3400\" \"39\u00d789=3471\"\nReplace-Text \"83\u00d771=5893\" \"85\u00d772=6120\"\nReplace-Text \"49\u00d728=1372\" \"21\u00d746=966\"\nReplace-Text \"86\u00d749=4214\" \"48\u00d760=2880\"\nReplace-Text \"78\u00d777=6006\" \"18\u00d755=990\"\nReplace-Text \"67\u00d778=5226\" \"60\u00d777=4620\"\nReplace-Text \"93\u00d756=5208\" \"44\u00d777=3388\"\nReplace-Text \"62\u00d741=2542\" \"63\u00d723=1449\"\nReplace-Text \"61\u00d750=3050\" \"68\u00d795=6460\"\nReplace-Text \"56\u00d799=5544\" \"75\u00d736=2700\"\nReplace-Text \"14\u00d785=1190\" \"33\u00d759=1947\"\nReplace-Text \"97\u00d712=1164\" \"22\u00d796=2112\"\nReplace-Text \"25\u00d728=700\" \"83\u00d717=1411\"\nReplace-Text \"48\u00d731=1488\" \"84\u00d734=2856\"\nReplace-Text \"14\u00d753=742\" \"14\u00d761=854\"\nReplace-Text \"62\u00d746=2852\" \"82\u00d752=4264\"\nReplace-Text \"91\u00d758=5278\" \"97\u00d716=1552\"\nReplace-Text \"17\u00d785=1445\" \"65\u00d762=4030\"\nReplace-Text \"57\u00d794=5358\" \"33\u00d742=1386\"\nReplace-Text \"19\u00d776=1444\" \"17\u00d763=1071\"\n"}
